# Auto-generated script applying the Behemoth_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 3843.7693  # H32: was 4078.6365
$ws.Cells.Item(32, 9).Value = 3442  # I32: was 3361.6667
$ws.Cells.Item(32, 10).Value = 4094.875  # J32: was 4939
$ws.Cells.Item(32, 11).Value = 3442  # K32: was 3361.6667
$ws.Cells.Item(32, 12).Value = 4094.875  # L32: was 4939
$ws.Cells.Item(32, 13).Value = -3116  # M32: was -3035.6667
$ws.Cells.Item(32, 14).Value = -4746.875  # N32: was -5591
$ws.Cells.Item(62, 8).Value = 5888  # H62: was 5886.8887
$ws.Cells.Item(62, 9).Value = 2997  # I62: was 2992
$ws.Cells.Item(62, 11).Value = 2997  # K62: was 2992
$ws.Cells.Item(62, 13).Value = -2373  # M62: was -2368
$ws.Cells.Item(65, 8).Value = 5888  # H65: was 5886.8887
$ws.Cells.Item(65, 9).Value = 2997  # I65: was 2992
$ws.Cells.Item(65, 11).Value = 14985  # K65: was 14960
$ws.Cells.Item(65, 13).Value = -11865  # M65: was -11840
$ws.Cells.Item(97, 8).Value = 1800  # H97: was 1650
$ws.Cells.Item(97, 10).Value = 1800  # J97: was 1650
$ws.Cells.Item(97, 12).Value = 5400  # L97: was 4950
$ws.Cells.Item(97, 14).Value = -6392  # N97: was -5942
$ws.Cells.Item(132, 8).Value = 782.90247  # H132: was 758.125
$ws.Cells.Item(132, 9).Value = 670.7632  # I132: was 671.7895
$ws.Cells.Item(132, 10).Value = 2203.3333  # J132: was 2398.5
$ws.Cells.Item(132, 11).Value = 2012.2896  # K132: was 2015.3685
$ws.Cells.Item(132, 12).Value = 6609.999899999999  # L132: was 7195.5
$ws.Cells.Item(132, 13).Value = 517.7103999999999  # M132: was 514.6315
$ws.Cells.Item(132, 14).Value = -11669.9999  # N132: was -12255.5
$ws.Cells.Item(138, 8).Value = 2080.465  # H138: was 2087.8572
$ws.Cells.Item(138, 10).Value = 2823.8215  # J138: was 2773.918
$ws.Cells.Item(138, 12).Value = 8471.4645  # L138: was 8321.754000000001
$ws.Cells.Item(138, 14).Value = -18751.4645  # N138: was -18601.754

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 80004  # H13: was 29868
$ws.Cells.Item(13, 9).Value = 0  # I13: was 1600
$ws.Cells.Item(13, 10).Value = 80004  # J13: was 44002
$ws.Cells.Item(13, 11).Value = 0  # K13: was 1600
$ws.Cells.Item(13, 12).Value = $null  # L13: was 44002
$ws.Cells.Item(13, 13).Value = 80004  # M13: was -1456
$ws.Cells.Item(13, 14).Value = -80292  # N13: was -44290
$ws.Cells.Item(32, 8).Value = 8068804.5  # H32: was 8337754.5
$ws.Cells.Item(32, 9).Value = 9805805  # I32: was 10206030
$ws.Cells.Item(32, 11).Value = 9805805  # K32: was 10206030
$ws.Cells.Item(32, 13).Value = -9805518  # M32: was -10205743
$ws.Cells.Item(46, 8).Value = 12749.5  # H46: was 15500
$ws.Cells.Item(46, 10).Value = 12749.5  # J46: was 15500
$ws.Cells.Item(46, 12).Value = 12749.5  # L46: was 15500
$ws.Cells.Item(46, 14).Value = -13387.5  # N46: was -16138
$ws.Cells.Item(122, 8).Value = 1991.6904  # H122: was 2022.7561
$ws.Cells.Item(122, 9).Value = 1268.5161  # I122: was 1286.8667
$ws.Cells.Item(122, 11).Value = 3805.5483  # K122: was 3860.6001
$ws.Cells.Item(122, 13).Value = -1355.5483  # M122: was -1410.6001
$ws.Cells.Item(138, 8).Value = 182164.67  # H138: was 182498
$ws.Cells.Item(138, 10).Value = 191597.6  # J138: was 191997.6
$ws.Cells.Item(138, 12).Value = 191597.6  # L138: was 191997.6
$ws.Cells.Item(138, 14).Value = -201877.6  # N138: was -202277.6

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2719.476  # H105: was 3133.25
$ws.Cells.Item(105, 9).Value = 2625.2856  # I105: was 5700
$ws.Cells.Item(105, 11).Value = 2625.2856  # K105: was 5700
$ws.Cells.Item(105, 13).Value = -878.2856000000002  # M105: was -3953

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 606371  # H31: was 622296.5600000001
$ws.Cells.Item(31, 9).Value = 9372.137000000001  # I31: was 9761.286
$ws.Cells.Item(31, 11).Value = 9372.137000000001  # K31: was 9761.286
$ws.Cells.Item(31, 13).Value = -9077.137000000001  # M31: was -9466.286
$ws.Cells.Item(34, 8).Value = 606371  # H34: was 622296.5600000001
$ws.Cells.Item(34, 9).Value = 9372.137000000001  # I34: was 9761.286
$ws.Cells.Item(34, 11).Value = 9372.137000000001  # K34: was 9761.286
$ws.Cells.Item(34, 13).Value = -9170.137000000001  # M34: was -9559.286
$ws.Cells.Item(42, 8).Value = 0  # H42: was 15056
$ws.Cells.Item(42, 9).Value = 0  # I42: was 15056
$ws.Cells.Item(42, 11).Value = 0  # K42: was 15056
$ws.Cells.Item(42, 13).Value = $null  # M42: was -14463
$ws.Cells.Item(58, 8).Value = 2074.25  # H58: was 2172
$ws.Cells.Item(58, 9).Value = 1630.909  # I58: was 1716.8572
$ws.Cells.Item(58, 11).Value = 1630.909  # K58: was 1716.8572
$ws.Cells.Item(58, 13).Value = -1427.909  # M58: was -1513.8572
$ws.Cells.Item(87, 8).Value = 68497.25  # H87: was 73139.86
$ws.Cells.Item(87, 10).Value = 52595.6  # J87: was 56744.75
$ws.Cells.Item(87, 12).Value = 52595.6  # L87: was 56744.75
$ws.Cells.Item(87, 14).Value = -54967.6  # N87: was -59116.75
$ws.Cells.Item(90, 8).Value = 68497.25  # H90: was 73139.86
$ws.Cells.Item(90, 10).Value = 52595.6  # J90: was 56744.75
$ws.Cells.Item(90, 12).Value = 157786.8  # L90: was 170234.25
$ws.Cells.Item(90, 14).Value = -169642.8  # N90: was -182090.25
$ws.Cells.Item(99, 8).Value = 3345.8235  # H99: was 3429.9375
$ws.Cells.Item(99, 9).Value = 2998.4443  # I99: was 3123.25
$ws.Cells.Item(99, 11).Value = 2998.4443  # K99: was 3123.25
$ws.Cells.Item(99, 13).Value = -1500.4443  # M99: was -1625.25
$ws.Cells.Item(126, 8).Value = 3345.8235  # H126: was 3429.9375
$ws.Cells.Item(126, 9).Value = 2998.4443  # I126: was 3123.25
$ws.Cells.Item(126, 11).Value = 8995.332900000001  # K126: was 9369.75
$ws.Cells.Item(126, 13).Value = -6525.332900000001  # M126: was -6899.75
$ws.Cells.Item(132, 8).Value = 2272.32  # H132: was 2309.1667
$ws.Cells.Item(132, 9).Value = 2278.652  # I132: was 2319.1365
$ws.Cells.Item(132, 11).Value = 6835.956  # K132: was 6957.4095
$ws.Cells.Item(132, 13).Value = -4305.956  # M132: was -4427.4095
$ws.Cells.Item(134, 8).Value = 2338.85  # H134: was 2373.3076
$ws.Cells.Item(134, 9).Value = 1871.5758  # I134: was 1898.9688
$ws.Cells.Item(134, 11).Value = 5614.7274  # K134: was 5696.9064
$ws.Cells.Item(134, 13).Value = -3079.7274  # M134: was -3161.9064
$ws.Cells.Item(136, 8).Value = 2074.25  # H136: was 2172
$ws.Cells.Item(136, 9).Value = 1630.909  # I136: was 1716.8572
$ws.Cells.Item(136, 11).Value = 4892.727000000001  # K136: was 5150.571599999999
$ws.Cells.Item(136, 13).Value = -2342.727000000001  # M136: was -2600.571599999999

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 37391304  # H4: was 43623172
$ws.Cells.Item(4, 9).Value = 18408434  # I4: was 23010516
$ws.Cells.Item(4, 11).Value = 55225302  # K4: was 69031548
$ws.Cells.Item(4, 13).Value = -55225190  # M4: was -69031436
$ws.Cells.Item(37, 8).Value = 84499  # H37: was 84461.75
$ws.Cells.Item(37, 10).Value = 84499  # J37: was 84461.75
$ws.Cells.Item(37, 12).Value = 253497  # L37: was 253385.25
$ws.Cells.Item(37, 14).Value = -253721  # N37: was -253609.25
$ws.Cells.Item(131, 8).Value = 6528.1387  # H131: was 7434.875
$ws.Cells.Item(131, 10).Value = 6687.971  # J131: was 7573.1914
$ws.Cells.Item(131, 12).Value = 20063.913  # L131: was 22719.5742
$ws.Cells.Item(131, 14).Value = -30143.913  # N131: was -32799.5742

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 169.4762  # H2: was 198.41176
$ws.Cells.Item(2, 9).Value = 66.333336  # I2: was 76.25
$ws.Cells.Item(2, 11).Value = 66.333336  # K2: was 76.25
$ws.Cells.Item(2, 13).Value = 46.666664  # M2: was 36.75
$ws.Cells.Item(20, 8).Value = 49554.4  # H20: was 19995.666
$ws.Cells.Item(20, 10).Value = 49554.4  # J20: was 19995.666
$ws.Cells.Item(20, 12).Value = 49554.4  # L20: was 19995.666
$ws.Cells.Item(20, 14).Value = -50044.4  # N20: was -20485.666
$ws.Cells.Item(24, 8).Value = 61110.5  # H24: was 0
$ws.Cells.Item(24, 10).Value = 61110.5  # J24: was 0
$ws.Cells.Item(24, 12).Value = 61110.5  # L24: was 0
$ws.Cells.Item(24, 14).Value = -61456.5  # N24: was None
$ws.Cells.Item(43, 8).Value = 9903.4  # H43: was 10999
$ws.Cells.Item(43, 9).Value = 10672.333  # I43: was 0
$ws.Cells.Item(43, 10).Value = 8750  # J43: was 10999
$ws.Cells.Item(43, 11).Value = 10672.333  # K43: was 0
$ws.Cells.Item(43, 12).Value = 8750  # L43: was 10999
$ws.Cells.Item(43, 13).Value = -10521.333  # M43: was None
$ws.Cells.Item(43, 14).Value = -9052  # N43: was -11301
$ws.Cells.Item(102, 8).Value = 3807.75  # H102: was 4494.9165
$ws.Cells.Item(102, 9).Value = 3431.2222  # I102: was 4278
$ws.Cells.Item(102, 10).Value = 4291.857  # J102: was 4798.6
$ws.Cells.Item(102, 11).Value = 3431.2222  # K102: was 4278
$ws.Cells.Item(102, 12).Value = 4291.857  # L102: was 4798.6
$ws.Cells.Item(102, 13).Value = -1809.2222  # M102: was -2656
$ws.Cells.Item(102, 14).Value = -7535.857  # N102: was -8042.6

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 2050.5  # H2: was 0
$ws.Cells.Item(2, 9).Value = 101  # I2: was 0
$ws.Cells.Item(2, 10).Value = 4000  # J2: was 0
$ws.Cells.Item(2, 11).Value = 101  # K2: was 0
$ws.Cells.Item(2, 12).Value = 4000  # L2: was 0
$ws.Cells.Item(2, 13).Value = 11  # M2: was None
$ws.Cells.Item(2, 14).Value = -4224  # N2: was None
$ws.Cells.Item(40, 8).Value = 3551.75  # H40: was 3268.4827
$ws.Cells.Item(40, 9).Value = 2048.7273  # I40: was 2128.7
$ws.Cells.Item(40, 10).Value = 5388.778  # J40: was 3868.3684
$ws.Cells.Item(40, 11).Value = 2048.7273  # K40: was 2128.7
$ws.Cells.Item(40, 12).Value = 5388.778  # L40: was 3868.3684
$ws.Cells.Item(40, 13).Value = -1912.7273  # M40: was -1992.7
$ws.Cells.Item(40, 14).Value = -5660.778  # N40: was -4140.368399999999
$ws.Cells.Item(55, 8).Value = 45455172  # H55: was 47619696
$ws.Cells.Item(55, 10).Value = 527.8570999999999  # J55: was 582.5
$ws.Cells.Item(55, 12).Value = 527.8570999999999  # L55: was 582.5
$ws.Cells.Item(55, 14).Value = -873.8570999999999  # N55: was -928.5
$ws.Cells.Item(61, 8).Value = 1560.5333  # H61: was 1572.2307
$ws.Cells.Item(61, 9).Value = 1525.25  # I61: was 1534.909
$ws.Cells.Item(61, 10).Value = 1701.6666  # J61: was 1777.5
$ws.Cells.Item(61, 11).Value = 1525.25  # K61: was 1534.909
$ws.Cells.Item(61, 12).Value = 1701.6666  # L61: was 1777.5
$ws.Cells.Item(61, 13).Value = -1323.25  # M61: was -1332.909
$ws.Cells.Item(61, 14).Value = -2105.6666  # N61: was -2181.5
$ws.Cells.Item(109, 8).Value = 100134.5  # H109: was 101640
$ws.Cells.Item(109, 10).Value = 100134.5  # J109: was 101640
$ws.Cells.Item(109, 12).Value = 100134.5  # L109: was 101640
$ws.Cells.Item(109, 14).Value = -102908.5  # N109: was -104414
$ws.Cells.Item(113, 8).Value = 1560.5333  # H113: was 1572.2307
$ws.Cells.Item(113, 9).Value = 1525.25  # I113: was 1534.909
$ws.Cells.Item(113, 10).Value = 1701.6666  # J113: was 1777.5
$ws.Cells.Item(113, 11).Value = 1525.25  # K113: was 1534.909
$ws.Cells.Item(113, 12).Value = 1701.6666  # L113: was 1777.5
$ws.Cells.Item(113, 13).Value = 644.75  # M113: was 635.0909999999999
$ws.Cells.Item(113, 14).Value = -6041.6666  # N113: was -6117.5
$ws.Cells.Item(123, 8).Value = 57985  # H123: was 57990
$ws.Cells.Item(123, 10).Value = 57985  # J123: was 57990
$ws.Cells.Item(123, 12).Value = 57985  # L123: was 57990
$ws.Cells.Item(123, 14).Value = -67785  # N123: was -67790
$ws.Cells.Item(132, 8).Value = 510631.2  # H132: was 537453.9
$ws.Cells.Item(132, 10).Value = 1669001.6  # J132: was 2002602
$ws.Cells.Item(132, 12).Value = 5007004.800000001  # L132: was 6007806
$ws.Cells.Item(132, 14).Value = -5012064.800000001  # N132: was -6012866

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 78340  # H24: was 83343.336
$ws.Cells.Item(24, 10).Value = 78340  # J24: was 83343.336
$ws.Cells.Item(24, 12).Value = 78340  # L24: was 83343.336
$ws.Cells.Item(24, 14).Value = -78800  # N24: was -83803.336
$ws.Cells.Item(31, 8).Value = 437259.5  # H31: was 437504.5
$ws.Cells.Item(31, 10).Value = 84519  # J31: was 85009
$ws.Cells.Item(31, 12).Value = 84519  # L31: was 85009
$ws.Cells.Item(31, 14).Value = -85215  # N31: was -85705
$ws.Cells.Item(51, 8).Value = 33681.5  # H51: was 40495
$ws.Cells.Item(51, 10).Value = 33681.5  # J51: was 40495
$ws.Cells.Item(51, 12).Value = 33681.5  # L51: was 40495
$ws.Cells.Item(51, 14).Value = -34701.5  # N51: was -41515
$ws.Cells.Item(68, 8).Value = 19100  # H68: was 28000
$ws.Cells.Item(68, 9).Value = 3650  # I68: was 0
$ws.Cells.Item(68, 10).Value = 50000  # J68: was 28000
$ws.Cells.Item(68, 11).Value = 3650  # K68: was 0
$ws.Cells.Item(68, 12).Value = 50000  # L68: was 28000
$ws.Cells.Item(68, 13).Value = -2839  # M68: was None
$ws.Cells.Item(68, 14).Value = -51622  # N68: was -29622
$ws.Cells.Item(70, 8).Value = 112990  # H70: was 88663
$ws.Cells.Item(70, 10).Value = 112990  # J70: was 88663
$ws.Cells.Item(70, 12).Value = 112990  # L70: was 88663
$ws.Cells.Item(70, 14).Value = -113620  # N70: was -89293
$ws.Cells.Item(71, 8).Value = 19100  # H71: was 28000
$ws.Cells.Item(71, 9).Value = 3650  # I71: was 0
$ws.Cells.Item(71, 10).Value = 50000  # J71: was 28000
$ws.Cells.Item(71, 11).Value = 10950  # K71: was 0
$ws.Cells.Item(71, 12).Value = 150000  # L71: was 84000
$ws.Cells.Item(71, 13).Value = -6894  # M71: was None
$ws.Cells.Item(71, 14).Value = -158112  # N71: was -92112
$ws.Cells.Item(73, 8).Value = 112990  # H73: was 88663
$ws.Cells.Item(73, 10).Value = 112990  # J73: was 88663
$ws.Cells.Item(73, 12).Value = 112990  # L73: was 88663
$ws.Cells.Item(73, 14).Value = -115174  # N73: was -90847
$ws.Cells.Item(103, 8).Value = 75277.75  # H103: was 61727.168
$ws.Cells.Item(103, 10).Value = 75277.75  # J103: was 61727.168
$ws.Cells.Item(103, 12).Value = 75277.75  # L103: was 61727.168
$ws.Cells.Item(103, 14).Value = -77621.75  # N103: was -64071.168
$ws.Cells.Item(109, 8).Value = 107990  # H109: was 108000
$ws.Cells.Item(109, 10).Value = 107990  # J109: was 108000
$ws.Cells.Item(109, 12).Value = 107990  # L109: was 108000
$ws.Cells.Item(109, 14).Value = -110764  # N109: was -110774
$ws.Cells.Item(122, 8).Value = 2660.2666  # H122: was 3040.8125
$ws.Cells.Item(122, 9).Value = 2454.4546  # I122: was 2979
$ws.Cells.Item(122, 11).Value = 7363.3638  # K122: was 8937
$ws.Cells.Item(122, 13).Value = -4913.3638  # M122: was -6487
$ws.Cells.Item(126, 8).Value = 1290.5714  # H126: was 1249.909
$ws.Cells.Item(126, 9).Value = 1290.5714  # I126: was 1224.9
$ws.Cells.Item(126, 10).Value = 0  # J126: was 1500
$ws.Cells.Item(126, 11).Value = 3871.7142  # K126: was 3674.7
$ws.Cells.Item(126, 12).Value = 0  # L126: was 4500
$ws.Cells.Item(126, 13).Value = $null  # M126: was -1204.7
$ws.Cells.Item(126, 14).Value = -1401.7142  # N126: was -9440
$ws.Cells.Item(140, 8).Value = 34317.375  # H140: was 33703.375
$ws.Cells.Item(140, 10).Value = 34317.375  # J140: was 33703.375
$ws.Cells.Item(140, 12).Value = 34317.375  # L140: was 33703.375
$ws.Cells.Item(140, 14).Value = -44677.375  # N140: was -44063.375
